# Clase 10 - corrige errores de tipeo / acentuacion detectados por el
# corrector ortografico y agrega signos de apertura de interrogacion.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1) "especifico" -> "específico" (parrafo sobre programar de forma especifica)
Replace-Text "tenemos que ser muy especifico al programar" "tenemos que ser muy específico al programar"

# 2) "mas" -> "más" y "rápida(" -> "rápida (" (antes de Javascript)
Replace-Text "escribir códigos de manera mas " "escribir códigos de manera más "
Replace-Text "rápida(" "rápida ("

# 3) "especificas" -> "específicas"
Replace-Text "tareas muy especificas y utilizar" "tareas muy específicas y utilizar"

# 4) "tambien" -> "también"
Replace-Text "desarrollar sino tambien a que hardware" "desarrollar sino también a que hardware"

# 5) "programacion" -> "programación"
Replace-Text "Lenguajes de programacion: " "Lenguajes de programación: "

# 6) Agrega el signo de apertura "¿" a las preguntas que no lo tenian
Replace-Text "Podríamos armar otro modelo?" "¿Podríamos armar otro modelo?"
Replace-Text "Por ahí uno que nos permite representar los objetos de la nave?" "¿Por ahí uno que nos permite representar los objetos de la nave?"

# 7) Agrega coma despues de "código"
Replace-Text "escribir código pero con otro set" "escribir código, pero con otro set"

# 8) "Alli" -> "Allí"
Replace-Text "Alli nació el lenguaje" "Allí nació el lenguaje"

# 9) Convierte el salto de linea (manual line break) antes de
# "¿Podemos elegir..." en un salto de parrafo real.
$findBreak = "ejecutable. " + [char]11 + "¿Podemos elegir"
$d.Content.Find.Execute($findBreak, $true, $false, $true, $false, $false, `
                         $true, 1, $false, "ejecutable. ^p¿Podemos elegir", 2) | Out-Null
